# FX - BOM update
#
# Several BOM line items were physically located/counted and are now fully
# in stock, so their on-hand quantity ("Have", column C) is set to match
# the required quantity ("Qty", column B). The old placeholder quantities
# that had been entered under "Proto" (D) or "Bought" (E) for those same
# rows are cleared out since the parts have now actually been counted as
# "Have".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FX-BOM")

# Row 2: Proto 1 -> Have 1
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 1

# Row 3: Proto 1 -> Have 1
$ws.Range("D3").ClearContents()
$ws.Range("C3").Value = 1

# Row 6: Bought 20 -> Have 3
$ws.Range("E6").ClearContents()
$ws.Range("C6").Value = 3

# Row 7: Have 2 & Bought 10 -> Have 4
$ws.Range("E7").ClearContents()
$ws.Range("C7").Value = 4

# Row 9: Have 0 -> cleared (quantity needed is 0, nothing to track)
$ws.Range("C9").ClearContents()

# Row 10: Proto 2 -> Have 2
$ws.Range("D10").ClearContents()
$ws.Range("C10").Value = 2

# Row 11: Proto 1 -> Have 1
$ws.Range("D11").ClearContents()
$ws.Range("C11").Value = 1

# Row 24: Bought 20 -> Have 4
$ws.Range("E24").ClearContents()
$ws.Range("C24").Value = 4

# Row 32: Proto 1 -> Have 1
$ws.Range("D32").ClearContents()
$ws.Range("C32").Value = 1

# Row 33: Proto 1 -> Have 1
$ws.Range("D33").ClearContents()
$ws.Range("C33").Value = 1

# Leave the selection on the last-edited cell, as the author did.
$ws.Range("C32").Select()
